$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Dinosaur', ['Token Creature $([char]0x2014) Dinosaur', 'Trample (This creature can deal excess combat damage to the player or planeswalker it$([char]0x2019)s attacking.)', '3/3'])"
$ws.Range("A3").Value = "('Dragon', ['Token Creature $([char]0x2014) Dragon', 'Flying', '5/5'])"

$ws.Range("A4:A14").EntireRow.Delete()
